$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 13.813694
$ws.Cells.Item(2, 8).Value = 41.441082
$ws.Cells.Item(2, 9).Value = 0.00564671359801415
$ws.Cells.Item(2, 10).Value = 0.005783813933563532
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.506715
$ws.Cells.Item(2, 14).Value = 1.520145
$ws.Cells.Item(2, 15).Value = 0.003122343715987576
$ws.Cells.Item(2, 16).Value = 0.003132472094339857
$ws.Cells.Item(2, 17).Value = 6.999605955210001
$ws.Cells.Item(2, 18).Value = 62.99645359689001
$ws.Cells.Item(2, 19).Value = [double]"1.763098071874108E-05"
$ws.Cells.Item(2, 20).Value = [double]"1.81176357457418E-05"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 13.813694
$ws.Cells.Item(3, 8).Value = 41.441082
$ws.Cells.Item(3, 9).Value = 0.00564671359801415
$ws.Cells.Item(3, 10).Value = 0.005783813933563532
$ws.Cells.Item(3, 13).Value = 88.13219433333332
$ws.Cells.Item(3, 14).Value = 264.396583
$ws.Cells.Item(3, 15).Value = 0.5430646480820168
$ws.Cells.Item(3, 16).Value = 0.5448262620252092
$ws.Cells.Item(3, 17).Value = 1217.4311640692
$ws.Cells.Item(3, 18).Value = 10956.8804766228
$ws.Cells.Item(3, 19).Value = 0.003066530532925493
$ws.Cells.Item(3, 20).Value = 0.003151173725672741
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 13.813694
$ws.Cells.Item(4, 8).Value = 41.441082
$ws.Cells.Item(4, 9).Value = 0.00564671359801415
$ws.Cells.Item(4, 10).Value = 0.005783813933563532
$ws.Cells.Item(4, 13).Value = 1.5741895
$ws.Cells.Item(4, 14).Value = 3.148379
$ws.Cells.Item(4, 15).Value = 0.009700049718478087
$ws.Cells.Item(4, 16).Value = 0.006487676741301404
$ws.Cells.Item(4, 17).Value = 21.745372051013
$ws.Cells.Item(4, 18).Value = 130.472232306078
$ws.Cells.Item(4, 19).Value = [double]"5.477340264674354E-05"
$ws.Cells.Item(4, 20).Value = [double]"3.752351513279511E-05"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 13.813694
$ws.Cells.Item(5, 8).Value = 41.441082
$ws.Cells.Item(5, 9).Value = 0.00564671359801415
$ws.Cells.Item(5, 10).Value = 0.005783813933563532
$ws.Cells.Item(5, 13).Value = 72.07364666666666
$ws.Cells.Item(5, 14).Value = 216.22094
$ws.Cells.Item(5, 15).Value = 0.4441129584835175
$ws.Cells.Item(5, 16).Value = 0.4455535891391496
$ws.Cells.Item(5, 17).Value = 995.6033005174532
$ws.Cells.Item(5, 18).Value = 8960.429704657079
$ws.Cells.Item(5, 19).Value = 0.002507778681723171
$ws.Cells.Item(5, 20).Value = 0.002576999057012254
$ws.Cells.Item(6, 7).Value = 60.71131066666667
$ws.Cells.Item(6, 9).Value = 0.02481735757971244
$ws.Cells.Item(6, 10).Value = 0.02541991480039814
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.506715
$ws.Cells.Item(6, 14).Value = 1.520145
$ws.Cells.Item(6, 15).Value = 0.003122343715987576
$ws.Cells.Item(6, 16).Value = 0.003132472094339857
$ws.Cells.Item(6, 17).Value = 30.76333178446
$ws.Cells.Item(6, 18).Value = 276.86998606014
$ws.Cells.Item(6, 19).Value = [double]"7.748832048643177E-05"
$ws.Cells.Item(6, 20).Value = [double]"7.962717375274388E-05"
$ws.Cells.Item(7, 7).Value = 60.71131066666667
$ws.Cells.Item(7, 9).Value = 0.02481735757971244
$ws.Cells.Item(7, 10).Value = 0.02541991480039814
$ws.Cells.Item(7, 13).Value = 88.13219433333332
$ws.Cells.Item(7, 14).Value = 264.396583
$ws.Cells.Item(7, 15).Value = 0.5430646480820168
$ws.Cells.Item(7, 16).Value = 0.5448262620252092
$ws.Cells.Item(7, 17).Value = 5350.621029906039
$ws.Cells.Item(7, 18).Value = 48155.58926915436
$ws.Cells.Item(7, 19).Value = 0.01347742956035211
$ws.Cells.Item(7, 20).Value = 0.01384943716170021
$ws.Cells.Item(8, 7).Value = 60.71131066666667
$ws.Cells.Item(8, 9).Value = 0.02481735757971244
$ws.Cells.Item(8, 10).Value = 0.02541991480039814
$ws.Cells.Item(8, 13).Value = 1.5741895
$ws.Cells.Item(8, 14).Value = 3.148379
$ws.Cells.Item(8, 15).Value = 0.009700049718478087
$ws.Cells.Item(8, 16).Value = 0.006487676741301404
$ws.Cells.Item(8, 17).Value = 95.57110778270467
$ws.Cells.Item(8, 18).Value = 573.426646696228
$ws.Cells.Item(8, 19).Value = 0.0002407296024044596
$ws.Cells.Item(8, 20).Value = 0.0001649161900164063
$ws.Cells.Item(9, 7).Value = 60.71131066666667
$ws.Cells.Item(9, 9).Value = 0.02481735757971244
$ws.Cells.Item(9, 10).Value = 0.02541991480039814
$ws.Cells.Item(9, 13).Value = 72.07364666666666
$ws.Cells.Item(9, 14).Value = 216.22094
$ws.Cells.Item(9, 15).Value = 0.4441129584835175
$ws.Cells.Item(9, 16).Value = 0.4455535891391496
$ws.Cells.Item(9, 17).Value = 4375.685553659564
$ws.Cells.Item(9, 18).Value = 39381.16998293608
$ws.Cells.Item(9, 19).Value = 0.01102171009646944
$ws.Cells.Item(9, 20).Value = 0.01132593427492878
$ws.Cells.Item(10, 7).Value = 521.5371296666666
$ws.Cells.Item(10, 8).Value = 1564.611389
$ws.Cells.Item(10, 9).Value = 0.2131921267372822
$ws.Cells.Item(10, 10).Value = 0.2183683609494171
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.506715
$ws.Cells.Item(10, 14).Value = 1.520145
$ws.Cells.Item(10, 15).Value = 0.003122343715987576
$ws.Cells.Item(10, 16).Value = 0.003132472094339857
$ws.Cells.Item(10, 17).Value = 264.270686659045
$ws.Cells.Item(10, 18).Value = 2378.436179931405
$ws.Cells.Item(10, 19).Value = 0.0006656590972161802
$ws.Cells.Item(10, 20).Value = 0.0006840327969607824
$ws.Cells.Item(11, 7).Value = 521.5371296666666
$ws.Cells.Item(11, 8).Value = 1564.611389
$ws.Cells.Item(11, 9).Value = 0.2131921267372822
$ws.Cells.Item(11, 10).Value = 0.2183683609494171
$ws.Cells.Item(11, 13).Value = 88.13219433333332
$ws.Cells.Item(11, 14).Value = 264.396583
$ws.Cells.Item(11, 15).Value = 0.5430646480820168
$ws.Cells.Item(11, 16).Value = 0.5448262620252092
$ws.Cells.Item(11, 17).Value = 45964.21166383152
$ws.Cells.Item(11, 18).Value = 413677.9049744837
$ws.Cells.Item(11, 19).Value = 0.1157771072804389
$ws.Cells.Item(11, 20).Value = 0.1189728178406426
$ws.Cells.Item(12, 7).Value = 521.5371296666666
$ws.Cells.Item(12, 8).Value = 1564.611389
$ws.Cells.Item(12, 9).Value = 0.2131921267372822
$ws.Cells.Item(12, 10).Value = 0.2183683609494171
$ws.Cells.Item(12, 13).Value = 1.5741895
$ws.Cells.Item(12, 14).Value = 3.148379
$ws.Cells.Item(12, 15).Value = 0.009700049718478087
$ws.Cells.Item(12, 16).Value = 0.006487676741301404
$ws.Cells.Item(12, 17).Value = 820.9982733814051
$ws.Cells.Item(12, 18).Value = 4925.989640288431
$ws.Cells.Item(12, 19).Value = 0.002067974228939719
$ws.Cells.Item(12, 20).Value = 0.001416703336367643
$ws.Cells.Item(13, 7).Value = 521.5371296666666
$ws.Cells.Item(13, 8).Value = 1564.611389
$ws.Cells.Item(13, 9).Value = 0.2131921267372822
$ws.Cells.Item(13, 10).Value = 0.2183683609494171
$ws.Cells.Item(13, 13).Value = 72.07364666666666
$ws.Cells.Item(13, 14).Value = 216.22094
$ws.Cells.Item(13, 15).Value = 0.4441129584835175
$ws.Cells.Item(13, 16).Value = 0.4455535891391496
$ws.Cells.Item(13, 17).Value = 37589.08280714284
$ws.Cells.Item(13, 18).Value = 338301.7452642856
$ws.Cells.Item(13, 19).Value = 0.09468138613068743
$ws.Cells.Item(13, 20).Value = 0.09729480697544608
$ws.Cells.Item(14, 7).Value = 173.964058
$ws.Cells.Item(14, 8).Value = 347.928116
$ws.Cells.Item(14, 9).Value = 0.07111242017336726
$ws.Cells.Item(14, 10).Value = 0.04855933745164542
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.506715
$ws.Cells.Item(14, 14).Value = 1.520145
$ws.Cells.Item(14, 15).Value = 0.003122343715987576
$ws.Cells.Item(14, 16).Value = 0.003132472094339857
$ws.Cells.Item(14, 17).Value = 88.15019764947
$ws.Cells.Item(14, 18).Value = 528.90118589682
$ws.Cells.Item(14, 19).Value = 0.0002220374182569814
$ws.Cells.Item(14, 20).Value = 0.0001521107694869116
$ws.Cells.Item(15, 7).Value = 173.964058
$ws.Cells.Item(15, 8).Value = 347.928116
$ws.Cells.Item(15, 9).Value = 0.07111242017336726
$ws.Cells.Item(15, 10).Value = 0.04855933745164542
$ws.Cells.Item(15, 13).Value = 88.13219433333332
$ws.Cells.Item(15, 14).Value = 264.396583
$ws.Cells.Item(15, 15).Value = 0.5430646480820168
$ws.Cells.Item(15, 16).Value = 0.5448262620252092
$ws.Cells.Item(15, 17).Value = 15331.83416667127
$ws.Cells.Item(15, 18).Value = 91991.00500002761
$ws.Cells.Item(15, 19).Value = 0.03861864143571021
$ws.Cells.Item(15, 20).Value = 0.02645640231020072
$ws.Cells.Item(16, 7).Value = 173.964058
$ws.Cells.Item(16, 8).Value = 347.928116
$ws.Cells.Item(16, 9).Value = 0.07111242017336726
$ws.Cells.Item(16, 10).Value = 0.04855933745164542
$ws.Cells.Item(16, 13).Value = 1.5741895
$ws.Cells.Item(16, 14).Value = 3.148379
$ws.Cells.Item(16, 15).Value = 0.009700049718478087
$ws.Cells.Item(16, 16).Value = 0.006487676741301404
$ws.Cells.Item(16, 17).Value = 273.852393480991
$ws.Cells.Item(16, 18).Value = 1095.409573923964
$ws.Cells.Item(16, 19).Value = 0.0006897940112829665
$ws.Cells.Item(16, 20).Value = 0.0003150372841580462
$ws.Cells.Item(17, 7).Value = 173.964058
$ws.Cells.Item(17, 8).Value = 347.928116
$ws.Cells.Item(17, 9).Value = 0.07111242017336726
$ws.Cells.Item(17, 10).Value = 0.04855933745164542
$ws.Cells.Item(17, 13).Value = 72.07364666666666
$ws.Cells.Item(17, 14).Value = 216.22094
$ws.Cells.Item(17, 15).Value = 0.4441129584835175
$ws.Cells.Item(17, 16).Value = 0.4455535891391496
$ws.Cells.Item(17, 17).Value = 12538.22404899151
$ws.Cells.Item(17, 18).Value = 75229.34429394903
$ws.Cells.Item(17, 19).Value = 0.0315819473081171
$ws.Cells.Item(17, 20).Value = 0.02163578708779974
$ws.Cells.Item(18, 7).Value = 1676.298339666667
$ws.Cells.Item(18, 8).Value = 5028.895019
$ws.Cells.Item(18, 9).Value = 0.685231381911624
$ws.Cells.Item(18, 10).Value = 0.7018685728649758
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 0.506715
$ws.Cells.Item(18, 14).Value = 1.520145
$ws.Cells.Item(18, 15).Value = 0.003122343715987576
$ws.Cells.Item(18, 16).Value = 0.003132472094339857
$ws.Cells.Item(18, 17).Value = 849.405513184195
$ws.Cells.Item(18, 18).Value = 7644.649618657755
$ws.Cells.Item(18, 19).Value = 0.002139527899309242
$ws.Cells.Item(18, 20).Value = 0.002198583718393677
$ws.Cells.Item(19, 7).Value = 1676.298339666667
$ws.Cells.Item(19, 8).Value = 5028.895019
$ws.Cells.Item(19, 9).Value = 0.685231381911624
$ws.Cells.Item(19, 10).Value = 0.7018685728649758
$ws.Cells.Item(19, 13).Value = 88.13219433333332
$ws.Cells.Item(19, 14).Value = 264.396583
$ws.Cells.Item(19, 15).Value = 0.5430646480820168
$ws.Cells.Item(19, 16).Value = 0.5448262620252092
$ws.Cells.Item(19, 17).Value = 147735.8510321466
$ws.Cells.Item(19, 18).Value = 1329622.65928932
$ws.Cells.Item(19, 19).Value = 0.3721249392725902
$ws.Cells.Item(19, 20).Value = 0.382396430986993
$ws.Cells.Item(20, 7).Value = 1676.298339666667
$ws.Cells.Item(20, 8).Value = 5028.895019
$ws.Cells.Item(20, 9).Value = 0.685231381911624
$ws.Cells.Item(20, 10).Value = 0.7018685728649758
$ws.Cells.Item(20, 13).Value = 1.5741895
$ws.Cells.Item(20, 14).Value = 3.148379
$ws.Cells.Item(20, 15).Value = 0.009700049718478087
$ws.Cells.Item(20, 16).Value = 0.006487676741301404
$ws.Cells.Item(20, 17).Value = 2638.8112451707
$ws.Cells.Item(20, 18).Value = 15832.8674710242
$ws.Cells.Item(20, 19).Value = 0.006646778473204199
$ws.Cells.Item(20, 20).Value = 0.004553496415626514
$ws.Cells.Item(21, 7).Value = 1676.298339666667
$ws.Cells.Item(21, 8).Value = 5028.895019
$ws.Cells.Item(21, 9).Value = 0.685231381911624
$ws.Cells.Item(21, 10).Value = 0.7018685728649758
$ws.Cells.Item(21, 13).Value = 72.07364666666666
$ws.Cells.Item(21, 14).Value = 216.22094
$ws.Cells.Item(21, 15).Value = 0.4441129584835175
$ws.Cells.Item(21, 16).Value = 0.4455535891391496
$ws.Cells.Item(21, 17).Value = 120816.9342410553
$ws.Cells.Item(21, 18).Value = 1087352.408169498
$ws.Cells.Item(21, 19).Value = 0.3043201362665204
$ws.Cells.Item(21, 20).Value = 0.3127200617439627
